$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '44.888.19'
$ws.Range('E2').Value = '  +0.12%  '
$ws.Range('D3').Value = '2.262.86'
$ws.Range('E3').Value = '  +0.50%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  -0.45%  '
$ws.Range('D5').Value = '''301.27'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = '''93.83'
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('D7').Value = '''0.564'
$ws.Range('E7').Value = '  -1.70%  '
$ws.Range('E8').Value = '  -0.49%  '
$ws.Range('D9').Value = '''0.508'
$ws.Range('E9').Value = '  -3.02%  '
$ws.Range('D10').Value = '''34.10'
$ws.Range('E10').Value = '  -4.13%  '
$ws.Range('D11').Value = '''0.0786'
$ws.Range('E11').Value = '  -2.48%  '
$ws.Range('D12').Value = '''7.17'
$ws.Range('E12').Value = '  -1.36%  '
$ws.Range('D13').Value = '''0.103'
$ws.Range('E13').Value = '  -1.04%  '
$ws.Range('D14').Value = '2.614.29'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').Value = '2.266.88'
$ws.Range('E15').Value = '  +0.60%  '
$ws.Range('D16').Value = '''13.66'
$ws.Range('E16').Value = '  +0.42%  '
$ws.Range('D17').Value = '''0.798'
$ws.Range('E17').Value = '  -5.39%  '
$ws.Range('D18').Value = '44.844.27'
$ws.Range('E18').Value = '  +0.68%  '
$ws.Range('D19').Value = '''12.79'
$ws.Range('E19').Value = '  +5.77%  '
$ws.Range('D20').Value = '0.0₃0918'
$ws.Range('E20').Value = '  -3.95%  '
$ws.Range('D21').Value = '''6.08'
$ws.Range('E21').Value = '  -3.97%  '
$ws.Range('D22').Value = '''65.28'
$ws.Range('E22').Value = '  -0.48%  '
$ws.Range('D23').Value = '''238.40'
$ws.Range('E23').Value = '  -0.17%  '
$ws.Range('D24').Value = '''2.88'
$ws.Range('E24').Value = '  -3.23%  '
$ws.Range('D25').Value = '''0.998'
$ws.Range('E25').Value = '  -0.39%  '
$ws.Range('D26').Value = '''1.91'
$ws.Range('E26').Value = '  -4.97%  '
$ws.Range('D27').Value = '''41.25'
$ws.Range('E27').Value = '  +10.17%  '
$ws.Range('E28').Value = '  -1.23%  '
$ws.Range('D29').Value = '''9.51'
$ws.Range('E29').Value = '  -3.91%  '
$ws.Range('D30').Value = '''19.52'
$ws.Range('E30').Value = '  -2.50%  '
$ws.Range('D31').Value = '''151.84'
$ws.Range('E31').Value = '  -0.56%  '
$ws.Range('D32').Value = '''5.55'
$ws.Range('E32').Value = '  -8.23%  '
$ws.Range('D33').Value = '''0.0788'
$ws.Range('E33').Value = '  -1.66%  '
$ws.Range('D34').Value = '''2.56'
$ws.Range('E34').Value = '  -2.85%  '
$ws.Range('D35').Value = '''2.92'
$ws.Range('E35').Value = '  -5.97%  '
$ws.Range('E36').Value = '  -2.08%  '
$ws.Range('E37').Value = '  -5.27%  '
$ws.Range('E38').Value = '  -6.10%  '
$ws.Range('B39').Value = 'VeChain'
$ws.Range('C39').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D39').Value = '''0.0306'
$ws.Range('E39').Value = '  +0.37%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''3.78'
$ws.Range('E40').Value = '  -2.07%  '
$ws.Range('D41').Value = '''3.23'
$ws.Range('E41').Value = '  -5.67%  '
$ws.Range('D42').Value = '''13.66'
$ws.Range('E42').Value = '  -8.68%  '
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  -0.80%  '
$ws.Range('D44').Value = '''1.94'
$ws.Range('E44').Value = '  +10.21%  '
$ws.Range('D45').Value = '1.761.10'
$ws.Range('E45').Value = '  -4.09%  '
$ws.Range('D46').Value = '''0.195'
$ws.Range('E46').Value = '  +1.07%  '
$ws.Range('D47').Value = '''69.96'
$ws.Range('E47').Value = '  -1.12%  '
$ws.Range('D48').Value = '''96.71'
$ws.Range('E48').Value = '  -2.93%  '
$ws.Range('D49').Value = '''75.25'
$ws.Range('E49').Value = '  -5.79%  '
$ws.Range('B50').Value = 'RocketPoolETH'
$ws.Range('C50').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D50').Value = '2.493.29'
$ws.Range('E50').Value = '  +0.69%  '
$ws.Range('D51').Value = '''52.98'
$ws.Range('E51').Value = '  -3.83%  '
